# Apply updated crypto price/volume figures (scraped refresh).
# Numeric-looking "Price" strings are written with a leading apostrophe so
# Excel keeps them as text (matching the original inline-string cells) instead
# of auto-converting them to numbers and dropping significant trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.747.76'
$ws.Range("E2").Value = '  -0.97%  '

$ws.Range("D3").Value = '3.493.24'
$ws.Range("E3").Value = '  -0.25%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").Value = "'600.41"
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").Value = "'147.49"
$ws.Range("E6").Value = '  -2.72%  '

$ws.Range("D7").Value = '3.492.64'
$ws.Range("E7").Value = '  -0.23%  '

$ws.Range("E8").Value = '  +0.06%  '

$ws.Range("D9").Value = "'0.478"
$ws.Range("E9").Value = '  -1.59%  '

$ws.Range("E10").Value = '  -0.84%  '

$ws.Range("E11").Value = '  +5.51%  '

$ws.Range("E12").Value = '  -2.50%  '

$ws.Range("D13").Value = "'0.0000212"
$ws.Range("E13").Value = '  -1.51%  '

$ws.Range("D14").Value = '4.085.45'
$ws.Range("E14").Value = '  +0.00%  '

$ws.Range("D15").Value = "'31.09"
$ws.Range("E15").Value = '  -4.31%  '

$ws.Range("D16").Value = '3.492.53'
$ws.Range("E16").Value = '  -0.17%  '

$ws.Range("D17").Value = '66.744.22'
$ws.Range("E17").Value = '  -0.85%  '

$ws.Range("E18").Value = '  +0.20%  '

$ws.Range("E19").Value = '  +6.47%  '

$ws.Range("D20").Value = "'6.34"
$ws.Range("E20").Value = '  -2.92%  '

$ws.Range("D21").Value = "'15.27"
$ws.Range("E21").Value = '  -1.55%  '

$ws.Range("D22").Value = "'432.94"
$ws.Range("E22").Value = '  -3.10%  '

$ws.Range("D23").Value = "'0.606"
$ws.Range("E23").Value = '  -4.06%  '

$ws.Range("D24").Value = "'79.64"
$ws.Range("E24").Value = '  +2.33%  '

$ws.Range("D25").Value = '3.633.16'
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("E26").Value = '  -0.14%  '

$ws.Range("E27").Value = '  -3.21%  '

$ws.Range("E28").Value = '  -7.44%  '

$ws.Range("D29").Value = "'9.83"
$ws.Range("E29").Value = '  -2.52%  '

$ws.Range("D30").Value = "'8.19"
$ws.Range("E30").Value = '  -7.25%  '

$ws.Range("E31").Value = '  -0.76%  '

$ws.Range("E32").Value = '  -2.78%  '

$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("E34").Value = '  -1.74%  '

$ws.Range("D35").Value = "'25.30"
$ws.Range("E35").Value = '  -1.60%  '

$ws.Range("D36").Value = '3.488.43'
$ws.Range("E36").Value = '  -0.01%  '

$ws.Range("D37").Value = "'5.87"
$ws.Range("E37").Value = '  -4.85%  '

$ws.Range("D38").Value = "'1.78"
$ws.Range("E38").Value = '  -4.64%  '

$ws.Range("D40").Value = "'7.96"
$ws.Range("E40").Value = '  -0.69%  '

$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("D42").Value = "'0.0890"
$ws.Range("E42").Value = '  -0.72%  '

$ws.Range("D43").Value = "'170.12"
$ws.Range("E43").Value = '  -2.47%  '

$ws.Range("E44").Value = '  -9.14%  '

$ws.Range("D45").Value = "'5.40"
$ws.Range("E45").Value = '  -1.06%  '

$ws.Range("D46").Value = "'0.895"
$ws.Range("E46").Value = '  +2.24%  '

$ws.Range("D47").Value = "'45.64"
$ws.Range("E47").Value = '  -3.03%  '

$ws.Range("D48").Value = "'28.21"
$ws.Range("E48").Value = '  -6.51%  '

$ws.Range("D49").Value = "'1.31"
$ws.Range("E49").Value = '  +0.65%  '

$ws.Range("D50").Value = "'7.44"
$ws.Range("E50").Value = '  -2.83%  '

$ws.Range("D51").Value = "'2.41"
$ws.Range("E51").Value = '  -4.23%  '
